# Eliminación de tildes y ñ's para el correcto funcionamiento de Linux
#
# Replace the element names that contain a tilde/diacritic or "ñ" with an
# accent-free version (the accented vowel becomes its uppercase plain
# counterpart, "ñ" becomes "N"), so the workbook behaves correctly when
# read back on Linux.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value  = "HidrOgeno"   # Hidrógeno
$ws.Range("B7").Value  = "NitrOgeno"   # Nitrógeno
$ws.Range("B8").Value  = "OxIgeno"     # Oxígeno
$ws.Range("B9").Value  = "FlUor"       # Flúor
$ws.Range("B10").Value = "NeOn"        # Neón
$ws.Range("B15").Value = "FOsforo"     # Fósforo
$ws.Range("B18").Value = "ArgOn"       # Argón
$ws.Range("B36").Value = "KriptOn"     # Kriptón
$ws.Range("B50").Value = "EstaNo"      # Estaño
$ws.Range("B54").Value = "XenOn"       # Xenón
$ws.Range("B73").Value = "TAntalo"     # Tántalo
$ws.Range("B86").Value = "RadOn"       # Radón

# Leave the view scrolled/selected near the end of the edited range, as in
# the source workbook.
$ws.Range("B119").Select()
